# Update "想去人数" (interested-count) figures in the "展览" and "全部类型"
# worksheets to the freshly scraped values.

$wb = $excel.ActiveWorkbook

# Map of cell address -> new value, applied identically to both sheets
# that hold this dataset ("展览" and "全部类型").
$updates = @{
    "F3"  = 393
    "F5"  = 8740
    "F7"  = 491
    "F9"  = 285
    "F14" = 365
    "F16" = 1227
    "F18" = 1123
    "F20" = 203
    "F21" = 2356
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
